$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# New benchmark rows 14-15 (VS2013 x64 4 core 8 thread i7-3770K, 1st run)
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "write VS2013 x64 4 core 8 thread Intel i7-3770K @ 3.50Ghz"
$ws.Range("B14").Formula = "=C14/D14"
$ws.Range("C14").Formula = "=H14/E14"
$ws.Range("D14").Value = 4
$ws.Range("E14").Value = 19897995
$ws.Range("H14").Value = 18289470

$ws.Range("A15").Value = "read VS2013 x64 4 core 8 thread Intel i7-3770K @ 3.50Ghz"
$ws.Range("B15").Formula = "=C15/D15"
$ws.Range("C15").Formula = "=H15/E15"
$ws.Range("D15").Value = 4
$ws.Range("E15").Value = 42025175
$ws.Range("H15").Value = 42351210

# ---------------------------------------------------------------------------
# Italic style for rows 31-32 (2 core results)
# ---------------------------------------------------------------------------
$ws.Range("A31:F32").Font.Italic = $true

# ---------------------------------------------------------------------------
# New benchmark rows 33-34 (VS2013 x64 4 core 8 thread i7-3770K, 2nd run)
# ---------------------------------------------------------------------------
$ws.Range("A33").Value = "write VS2013 x64 4 core 8 thread Intel i7-3770K @ 3.50Ghz"
$ws.Range("B33").Formula = "=C33/D33"
$ws.Range("C33").Formula = "=H33/E33"
$ws.Range("D33").Value = 4
$ws.Range("E33").Value = 13458759
$ws.Range("H33").Value = 53613962

$ws.Range("A34").Value = "read VS2013 x64 4 core 8 thread Intel i7-3770K @ 3.50Ghz"
$ws.Range("B34").Formula = "=C34/D34"
$ws.Range("C34").Formula = "=H34/E34"
$ws.Range("D34").Value = 4
$ws.Range("E34").Value = 24216025
$ws.Range("H34").Value = 98416354

# ---------------------------------------------------------------------------
# Update the "Scaling" chart source formulas so they average in the new data
# ---------------------------------------------------------------------------
$ws.Range("D41").Formula = "=(C23+C33)/2"
$ws.Range("D42").Formula = "=C34"

# ---------------------------------------------------------------------------
# New data block (rows 46-48) feeding the second chart
# ---------------------------------------------------------------------------
$ws.Range("B46").Value = 1
$ws.Range("C46").Value = 2
$ws.Range("D46").Value = 4
$ws.Range("E46").Value = 8

$ws.Range("A47").Value = "GCC 4.8 write over VS2013 write"
$ws.Range("B47").Formula = "=E29/E33"
$ws.Range("C47").Formula = "=E31/E25/(3.5/2.53)"
$ws.Range("D47").Formula = "=H29/H33"
$ws.Range("C47").Font.Italic = $true

$ws.Range("A48").Value = "GCC 4.8 write over VS2013 read"
$ws.Range("B48").Formula = "=E30/E34"
$ws.Range("C48").Formula = "=E32/E26/(3.5/2.53)"
$ws.Range("D48").Formula = "=H30/H34"
$ws.Range("C48").Font.Italic = $true

$excel.Calculate()

# ---------------------------------------------------------------------------
# Reposition the existing "Scaling" chart (data/series left untouched so the
# chart part itself is preserved as-is)
# ---------------------------------------------------------------------------
$co1 = $ws.ChartObjects().Item(1)
$co1.Left = 949.6249212598425
$co1.Top = 253.5
$co1.Width = 568.625
$co1.Height = 337.5

# ---------------------------------------------------------------------------
# New chart: "By how many times GCC 4.8 is faster than VS2013 ..."
# ---------------------------------------------------------------------------
$co2 = $ws.ChartObjects().Add(948.1249212598425, 597.0, 570.125, 301.5)
$chart2 = $co2.Chart
$chart2.ChartType = 74

$ser2_1 = $chart2.SeriesCollection().NewSeries()
$ser2_1.Name = "=Sheet1!`$A`$47"
$ser2_1.XValues = $ws.Range("B46:D46")
$ser2_1.Values = $ws.Range("B47:D47")

$ser2_2 = $chart2.SeriesCollection().NewSeries()
$ser2_2.Name = "=Sheet1!`$A`$48"
$ser2_2.XValues = $ws.Range("B46:D46")
$ser2_2.Values = $ws.Range("B48:D48")

$chart2.HasTitle = $true
$chart2.ChartTitle.Text = "By how many times GCC 4.8 is faster than VS2013 for proposed boost::concurrent_unordered_map"

$ax2_1 = $chart2.Axes(1)
$ax2_1.HasTitle = $true
$ax2_1.AxisTitle.Text = "CPU cores"
$ax2_1.MinimumScale = 1
$ax2_1.MaximumScale = 8

$ax2_2 = $chart2.Axes(2)
$ax2_2.HasTitle = $true
$ax2_2.AxisTitle.Text = "How many times faster"
$ax2_2.MinimumScale = 1
$ax2_2.HasMajorGridlines = $true

$chart2.HasLegend = $true
$chart2.Legend.Position = -4107

$co2.Left = 948.1249212598425
$co2.Top = 597.0
$co2.Width = 570.125
$co2.Height = 301.5

# ---------------------------------------------------------------------------
# Restore the selection on the sheet
# ---------------------------------------------------------------------------
$ws.Range("H48").Select()
